$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: "13" -> "9", F2: "1" -> "0"
$ws.Range("D2").Formula = '="9"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("F2").Formula = '="0"'
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)

# D4: "9" -> "13", F4: "0" -> "1"
$ws.Range("D4").Formula = '="13"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)

$ws.Range("F4").Formula = '="1"'
$ws.Range("F4").Copy()
$ws.Range("F4").PasteSpecial(-4163)
